$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: remove spaces -> underscores
$ws.Range("F1").Value = "EMISSION_CATEGORY"
$ws.Range("G1").Value = "MARKET_CATEGORY"
$ws.Range("B1").Value = "PRODUCT_BARCODE"

# Change the selection to C2
$ws.Range("C2").Select()
